$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($sheet, $rowNum, $values)
    $n = $values.Length
    $arr = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0, $i] = $values[$i]
    }
    $startCell = $sheet.Cells.Item($rowNum, 1)
    $endCell = $sheet.Cells.Item($rowNum, $n)
    $rng = $sheet.Range($startCell, $endCell)
    $rng.Value = $arr
}

# Insert two new rows at position 11 for the Swiss Super League matches,
# shifting the existing rows 11-14 down to 13-16.
$ws.Rows("11:12").Insert()

Set-RowValues $ws 2 @("Mexican Liga MX","'2026-01-14","'00:10:00","FC Juarez","Guadalajara",3.9,4.4,2.08,2.16,3.35,3.55,1.01,1.08,3.35,1.37,1.82,2.08,1.32,3.8,1.85,2.02,1.86,1.3,13,9,16,970,17,9.2,13,29,36,21,22,970,1000,70,65,140,65,22)
Set-RowValues $ws 3 @("Saudi Professional League","'2026-01-14","'11:45:00","Al-Quadisiya (KSA)","Al-Feiha",1.44,1.6,7.4,10.5,4.2,5.2,1.01,1.06,2.02,1.28,2.02,1.79,1.31,2.78,1.86,1.7,1.1,2.66,1000,38,1000,1000,12,15.5,46,1000,13,15,36,1000,19,23,1000,1000,1000,1000)
Set-RowValues $ws 4 @("German Bundesliga","'2026-01-14","'14:30:00","Wolfsburg","St Pauli",1.98,1.99,4.2,4.4,3.8,3.85,1.36,1.06,4,1.31,2,1.96,1.39,3.4,1.82,2.14,1.3,2,15,16,32,90,9.6,8.2,16.5,55,12,10,18.5,60,23,20,36,110,13.5,50)
Set-RowValues $ws 5 @("Italian Serie A","'2026-01-14","'14:30:00","Napoli","Parma",1.34,1.36,12,13.5,5.6,5.7,1.38,1.06,3.9,1.32,1.99,1.96,1.37,3.5,2.46,1.63,1.08,3.8,16,34,120,740,6.8,12,50,300,6.8,11,40,1000,9.8,16,55,290,6.8,500)
Set-RowValues $ws 6 @("Saudi Professional League","'2026-01-14","'14:30:00","Al-Shabab (KSA)","NEOM Sports Club",2.7,3.15,2.52,2.82,3.4,3.95,1.01,1.06,1.98,1.3,1.98,1.82,1.32,2.74,1.56,2,1.54,1.46,1000,16.5,25,1000,18,11,17,38,29,18.5,23,1000,1000,44,1000,1000,1000,1000)
Set-RowValues $ws 7 @("Saudi Professional League","'2026-01-14","'14:30:00","Al Ahli","Al-Taawoun Buraidah",1.57,1.71,5.4,7.4,4.2,5,1.01,1.05,4.1,1.25,2.16,1.64,1.43,2.84,1.81,1.99,1.15,2.4,23,26,1000,1000,11.5,12.5,27,1000,12.5,12.5,26,1000,19,21,42,1000,10,1000)
Set-RowValues $ws 8 @("German Bundesliga","'2026-01-14","'16:30:00","FC Koln","Bayern Munich",12,12.5,1.26,1.27,7.4,7.8,1.14,1.01,10.5,1.09,4.2,1.29,2.3,1.72,1.66,2.44,4.7,1.08,60,19.5,13,12.5,70,19.5,12,12,130,46,25,24,500,140,85,90,85,2.86)
Set-RowValues $ws 9 @("German Bundesliga","'2026-01-14","'16:30:00","Hoffenheim","Mgladbach",1.92,1.94,4.1,4.3,4.1,4.2,1.28,1.04,5.3,1.21,2.48,1.64,1.59,2.58,1.63,2.48,1.31,2.06,23,21,34,75,12.5,10,17,42,14.5,11,16.5,55,22,18,28,80,9,34)
Set-RowValues $ws 10 @("German Bundesliga","'2026-01-14","'16:30:00","RB Leipzig","Freiburg",1.79,1.81,4.8,4.9,4.2,4.3,1.29,1.04,5.3,1.21,2.46,1.65,1.59,2.62,1.66,2.42,1.25,2.24,22,22,40,110,12,9.6,19,50,12.5,10,16.5,50,19.5,17,28,75,8.2,42)
Set-RowValues $ws 11 @("Swiss Super League","'2026-01-14","'16:30:00","Servette","Lausanne",2.14,2.3,3.2,3.55,3.7,4.1,1.01,1.04,2.42,1.19,2.4,1.61,1.5,2.22,1.38,2.02,1.4,1.79,1000,26,40,85,20,13.5,21,48,24,17,21,50,40,30,40,85,16,1000)
Set-RowValues $ws 12 @("Swiss Super League","'2026-01-14","'16:30:00","Sion","Winterthur",1.48,1.52,7,8,5,5.2,1.01,1.04,2.46,1.18,2.46,1.62,1.5,2.22,1.51,1.79,1.14,2.94,1000,46,95,1000,15.5,17,40,1000,15.5,15,30,100,20,21,40,1000,8.2,1000)
Set-RowValues $ws 13 @("Italian Serie A","'2026-01-14","'16:45:00","Inter","Lecce",1.15,1.16,27,30,10,10.5,1.26,1.03,5.9,1.18,2.7,1.57,1.67,2.4,2.82,1.51,1.03,7.2,30,75,400,1000,9.4,23,110,1000,7,14.5,70,490,7.4,16,70,500,3.35,1000)
Set-RowValues $ws 14 @("Scottish Premiership","'2026-01-14","'16:45:00","Hearts","St Mirren",1.45,1.46,9.2,9.4,4.8,5,1.37,1.06,3.85,1.29,2,1.83,1.38,3.15,2.1,1.76,1.11,3.15,17.5,32,80,350,7.8,11,36,160,8,10.5,30,170,12,16,42,190,7.6,250)
Set-RowValues $ws 15 @("Scottish Premiership","'2026-01-14","'16:45:00","Falkirk","Celtic",6.8,7.6,1.49,1.53,4.7,5.2,1.28,1.04,5,1.2,2.38,1.6,1.55,2.52,1.79,2.06,2.88,1.15,26,12,10,15,32,12,11.5,16,65,27,23,30,200,100,85,110,95,6.6)
Set-RowValues $ws 16 @("Mexican Liga MX","'2026-01-14","'20:00:00","Cruz Azul","Atlas",1.31,1.5,7.4,16,4.5,10,1.01,1.01,2.16,1.2,2.16,1.56,1.52,2.32,1.9,1.94,1.06,2.96,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000,1000)
